$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Paragraph 4 ("3.3. Press X to go to the main page") is renumbered to "3.4."
# once the new "3.3. Click on Save" paragraph is inserted before it.
$para4 = $tr.Paragraphs(4)
$para4.Text = "3.4. Press X to go to the main page"

# Paragraph 3 ("3.2. Click on Save") becomes the new ChainId instruction.
$para3 = $tr.Paragraphs(3)
$para3.Text = "3.2. In the “ChainId” box, insert:"

# Insert the (now renumbered) "Click on Save" step right before the
# "Press X" paragraph.
$paraSave = $para4.InsertBefore("3.3. Click on Save`r")

# Insert the purple ChainId value paragraph right before "3.3. Click on Save",
# matching the styling used by the RPC URL value paragraph above it.
$paraChainId = $paraSave.InsertBefore("  `r")
$paraChainId.Font.Name = "Helvetica Neue Light"
$paraChainId.Font.Size = 24
$paraChainId.Font.Color.RGB = 9576852

$chainIdValue = $paraChainId.InsertAfter("13777222009")
$chainIdValue.Font.Name = "Helvetica Neue Light"
$chainIdValue.Font.Size = 24
$chainIdValue.Font.Color.RGB = 9576852

# The shape grew by two extra lines of text; the author nudged its
# position/size to account for the additional content.
$shp.Top = 218.0976377952756
$shp.Width = 873.8362204724409
$shp.Height = 494.844101488189
